$wb = $excel.ActiveWorkbook

# Sheet ALC, row 5
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 113.1875
$ws.Range("I5").Value = 29.428572
$ws.Range("J5").Value = 699.5
$ws.Range("K5").Value = 29.428572
$ws.Range("L5").Value = 699.5
$ws.Range("M5").Value = 85.571428
$ws.Range("N5").Value = -929.5

# Sheet ALC, row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1096.7
$ws.Range("J112").Value = 1141.0781
$ws.Range("L112").Value = 3423.2343
$ws.Range("N112").Value = -5639.2343

# Sheet ALC, row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1557.6571
$ws.Range("I137").Value = 1202.8846
$ws.Range("J137").Value = 2582.5557
$ws.Range("K137").Value = 3608.6538
$ws.Range("L137").Value = 7747.6671
$ws.Range("M137").Value = -1058.6538
$ws.Range("N137").Value = -12847.6671

# Sheet ALC, row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4875.8237
$ws.Range("I138").Value = 911.7353000000001
$ws.Range("J138").Value = 12804
$ws.Range("K138").Value = 2735.2059
$ws.Range("L138").Value = 38412
$ws.Range("N138").Value = -48692
$ws.Range("M138").Value = 2404.7941

# Sheet ALC, row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 1030.1428
$ws.Range("I141").Value = 1030.1428
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 3090.4284
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 2089.5716
$ws.Range("N141").ClearContents()

# Sheet ARM, row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5053.7407
$ws.Range("I32").Value = 3703.7812
$ws.Range("K32").Value = 3703.7812
$ws.Range("M32").Value = -3416.7812

# Sheet ARM, row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 10870868
$ws.Range("I74").Value = 1211.5161
$ws.Range("J74").Value = 33334826
$ws.Range("K74").Value = 1211.5161
$ws.Range("L74").Value = 33334826
$ws.Range("M74").Value = -337.5161000000001
$ws.Range("N74").Value = -33336574

# Sheet ARM, row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 10870868
$ws.Range("I77").Value = 1211.5161
$ws.Range("J77").Value = 33334826
$ws.Range("K77").Value = 6057.5805
$ws.Range("L77").Value = 166674130
$ws.Range("M77").Value = -1689.5805
$ws.Range("N77").Value = -166682866

# Sheet ARM, row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1765654.4
$ws.Range("I102").Value = 1951302.1
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 1951302.1
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = -1949680.1
$ws.Range("N102").Value = -5244

# Sheet BSM, row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 37038736
$ws.Range("I86").Value = 47620650
$ws.Range("J86").Value = 2050
$ws.Range("K86").Value = 47620650
$ws.Range("L86").Value = 2050
$ws.Range("M86").Value = -47619527
$ws.Range("N86").Value = -4296

# Sheet BSM, row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 37038736
$ws.Range("I89").Value = 47620650
$ws.Range("J89").Value = 2050
$ws.Range("K89").Value = 238103250
$ws.Range("L89").Value = 10250
$ws.Range("M89").Value = -238097634
$ws.Range("N89").Value = -21482

# Sheet BSM, row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 62501108
$ws.Range("I99").Value = 83334350
$ws.Range("J99").Value = 1375
$ws.Range("K99").Value = 83334350
$ws.Range("L99").Value = 1375
$ws.Range("M99").Value = -83332852
$ws.Range("N99").Value = -4371

# Sheet BSM, row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3373.0679
$ws.Range("I134").Value = 4040.8718
$ws.Range("J134").Value = 2070.85
$ws.Range("K134").Value = 12122.6154
$ws.Range("L134").Value = 6212.549999999999
$ws.Range("M134").Value = -9587.615399999999
$ws.Range("N134").Value = -11282.55

# Sheet CRP, row 16
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4274808.5
$ws.Range("I16").Value = 12821308
$ws.Range("J16").Value = 1558.4166
$ws.Range("K16").Value = 12821308
$ws.Range("L16").Value = 1558.4166
$ws.Range("M16").Value = -12821021
$ws.Range("N16").Value = -2132.4166

# Sheet CRP, row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10422302
$ws.Range("I31").Value = 2193.5417
$ws.Range("J31").Value = 20842410
$ws.Range("K31").Value = 2193.5417
$ws.Range("L31").Value = 20842410
$ws.Range("M31").Value = -1898.5417
$ws.Range("N31").Value = -20843000

# Sheet CRP, row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 10422302
$ws.Range("I34").Value = 2193.5417
$ws.Range("J34").Value = 20842410
$ws.Range("K34").Value = 2193.5417
$ws.Range("L34").Value = 20842410
$ws.Range("M34").Value = -1991.5417
$ws.Range("N34").Value = -20842814

# Sheet CRP, row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 5210074.5
$ws.Range("I58").Value = 9805003
$ws.Range("J58").Value = 2488.7334
$ws.Range("K58").Value = 9805003
$ws.Range("L58").Value = 2488.7334
$ws.Range("M58").Value = -9804800
$ws.Range("N58").Value = -2894.7334

# Sheet CRP, row 113
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 4274808.5
$ws.Range("I113").Value = 12821308
$ws.Range("J113").Value = 1558.4166
$ws.Range("K113").Value = 12821308
$ws.Range("L113").Value = 1558.4166
$ws.Range("M113").Value = -12819138
$ws.Range("N113").Value = -5898.4166

# Sheet CRP, row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 6062575.5
$ws.Range("I134").Value = 8132213
$ws.Range("J134").Value = 1495
$ws.Range("K134").Value = 24396639
$ws.Range("L134").Value = 4485
$ws.Range("M134").Value = -24394104
$ws.Range("N134").Value = -9555

# Sheet CRP, row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 5210074.5
$ws.Range("I136").Value = 9805003
$ws.Range("J136").Value = 2488.7334
$ws.Range("K136").Value = 29415009
$ws.Range("L136").Value = 7466.2002
$ws.Range("M136").Value = -29412459
$ws.Range("N136").Value = -12566.2002

# Sheet LTW, row 100
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1690.25
$ws.Range("I100").Value = 1680.5
$ws.Range("J100").Value = 1700
$ws.Range("K100").Value = 1680.5
$ws.Range("L100").Value = 1700
$ws.Range("M100").Value = -1139.5
$ws.Range("N100").Value = -2782

# Sheet LTW, row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 9968731
$ws.Range("I132").Value = 12392517
$ws.Range("K132").Value = 37177551
$ws.Range("M132").Value = -37175021

# Sheet LTW, row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 7061.977
$ws.Range("I136").Value = 5381.3
$ws.Range("J136").Value = 10663.429
$ws.Range("K136").Value = 16143.9
$ws.Range("L136").Value = 31990.287
$ws.Range("M136").Value = -13593.9
$ws.Range("N136").Value = -37090.287

# Sheet WVR, row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1141.4507
$ws.Range("I132").Value = 783.35297
$ws.Range("J132").Value = 2054.6
$ws.Range("K132").Value = 2350.05891
$ws.Range("L132").Value = 6163.799999999999
$ws.Range("M132").Value = 179.9410899999998
$ws.Range("N132").Value = -11223.8

# Sheet WVR, row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 4170261.8
$ws.Range("I136").Value = 5947.2354
$ws.Range("J136").Value = 7248233.5
$ws.Range("K136").Value = 17841.7062
$ws.Range("L136").Value = 21744700.5
$ws.Range("M136").Value = -15291.7062
$ws.Range("N136").Value = -21749800.5
